$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# =========================================================================
# 1) Rows 137 and 138 swap their match data (id/date/league stay put)
# =========================================================================
$ws.Range("B137").Value = 7499442
$ws.Range("F137").Value = "Guarani Asuncion"
$ws.Range("G137").Value = "Sportivo Trinidense"
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 2
$ws.Range("J137").Value = "A"
$ws.Range("K137").Value = 2
$ws.Range("L137").Value = 3.4
$ws.Range("M137").Value = 3.3
$ws.Range("N137").Value = 1.909
$ws.Range("O137").Value = 3.4
$ws.Range("P137").Value = 3.6
$ws.Range("Q137").Value = -0.5
$ws.Range("R137").Value = 1.975
$ws.Range("S137").Value = 1.825
$ws.Range("T137").Value = 2.5
$ws.Range("U137").Value = 1.9
$ws.Range("V137").Value = 1.9
$ws.Range("W137").Value = -1
$ws.Range("X137").Value = -1
$ws.Range("Y137").Value = 2.6
$ws.Range("Z137").Value = -1
$ws.Range("AA137").Value = 0.825
$ws.Range("AB137").Value = -1
$ws.Range("AC137").Value = 0.8999999999999999

$ws.Range("B138").Value = 7499440
$ws.Range("F138").Value = "Olimpia Asuncion"
$ws.Range("G138").Value = "Libertad Asuncion"
$ws.Range("H138").Value = 1
$ws.Range("I138").Value = 3
$ws.Range("J138").Value = "A"
$ws.Range("K138").Value = 2.7
$ws.Range("L138").Value = 3.4
$ws.Range("M138").Value = 2.3
$ws.Range("N138").Value = 2.8
$ws.Range("O138").Value = 3.25
$ws.Range("P138").Value = 2.3
$ws.Range("Q138").Value = 0.25
$ws.Range("R138").Value = 1.75
$ws.Range("S138").Value = 2.05
$ws.Range("T138").Value = 2.25
$ws.Range("U138").Value = 1.85
$ws.Range("V138").Value = 1.95
$ws.Range("W138").Value = -1
$ws.Range("X138").Value = -1
$ws.Range("Y138").Value = 1.3
$ws.Range("Z138").Value = -1
$ws.Range("AA138").Value = 1.05
$ws.Range("AB138").Value = 0.8500000000000001
$ws.Range("AC138").Value = -1

# =========================================================================
# 2) Rows 143 -> 144 -> 145 -> 143 cyclic rotation of match data
# =========================================================================
$ws.Range("B144").Value = 7493311
$ws.Range("F144").Value = "General Caballero JLM"
$ws.Range("G144").Value = "Olimpia Asuncion"
$ws.Range("H144").Value = 0
$ws.Range("I144").Value = 1
$ws.Range("J144").Value = "A"
$ws.Range("K144").Value = 3.4
$ws.Range("L144").Value = 3.3
$ws.Range("M144").Value = 2
$ws.Range("N144").Value = 3.2
$ws.Range("O144").Value = 3.25
$ws.Range("P144").Value = 2.1
$ws.Range("Q144").Value = 0.25
$ws.Range("R144").Value = 1.95
$ws.Range("S144").Value = 1.85
$ws.Range("T144").Value = 2.25
$ws.Range("U144").Value = 1.775
$ws.Range("V144").Value = 2.025
$ws.Range("W144").Value = -1
$ws.Range("X144").Value = -1
$ws.Range("Y144").Value = 1.1
$ws.Range("Z144").Value = -1
$ws.Range("AA144").Value = 0.8500000000000001
$ws.Range("AB144").Value = -1
$ws.Range("AC144").Value = 1.025

$ws.Range("B145").Value = 7493312
$ws.Range("F145").Value = "Cerro Porteno"
$ws.Range("G145").Value = "Guarani Asuncion"
$ws.Range("H145").Value = 4
$ws.Range("I145").Value = 0
$ws.Range("J145").Value = "H"
$ws.Range("K145").Value = 1.7
$ws.Range("L145").Value = 3.6
$ws.Range("M145").Value = 4.333
$ws.Range("N145").Value = 1.727
$ws.Range("O145").Value = 3.75
$ws.Range("P145").Value = 4.2
$ws.Range("Q145").Value = -0.5
$ws.Range("R145").Value = 1.8
$ws.Range("S145").Value = 2
$ws.Range("T145").Value = 2.75
$ws.Range("U145").Value = 1.875
$ws.Range("V145").Value = 1.925
$ws.Range("W145").Value = 0.7270000000000001
$ws.Range("X145").Value = -1
$ws.Range("Y145").Value = -1
$ws.Range("Z145").Value = 0.8
$ws.Range("AA145").Value = -1
$ws.Range("AB145").Value = 0.875
$ws.Range("AC145").Value = -1

$ws.Range("B143").Value = 7493433
$ws.Range("F143").Value = "Sportivo Luqueno"
$ws.Range("G143").Value = "Nacional Asuncion"
$ws.Range("H143").Value = 1
$ws.Range("I143").Value = 1
$ws.Range("J143").Value = "D"
$ws.Range("K143").Value = 2.75
$ws.Range("L143").Value = 3.2
$ws.Range("M143").Value = 2.4
$ws.Range("N143").Value = 2.75
$ws.Range("O143").Value = 3.1
$ws.Range("P143").Value = 2.45
$ws.Range("Q143").Value = 0.25
$ws.Range("R143").Value = 1.75
$ws.Range("S143").Value = 2.05
$ws.Range("T143").Value = 2.25
$ws.Range("U143").Value = 2
$ws.Range("V143").Value = 1.8
$ws.Range("W143").Value = -1
$ws.Range("X143").Value = 2.1
$ws.Range("Y143").Value = -1
$ws.Range("Z143").Value = 0.375
$ws.Range("AA143").Value = -0.5
$ws.Range("AB143").Value = -0.5
$ws.Range("AC143").Value = 0.4

# =========================================================================
# 3) Updated odds for upcoming fixtures (rows 210, 211, 213)
# =========================================================================
$ws.Range("N210").Value = 2.375
$ws.Range("P210").Value = 2.9
$ws.Range("Q210").Value = -0.25
$ws.Range("R210").Value = 2.05
$ws.Range("S210").Value = 1.75

$ws.Range("N211").Value = 2.45
$ws.Range("P211").Value = 2.8
$ws.Range("R211").Value = 1.825
$ws.Range("S211").Value = 1.975

$ws.Range("R213").Value = 2
$ws.Range("S213").Value = 1.8

# =========================================================================
# 4) Two new upcoming fixtures appended as rows 214 and 215
#    (copy number/date formatting from row 213, the last existing row)
# =========================================================================
$ws.Range("A213").Copy() | Out-Null
$ws.Range("A214").PasteSpecial(-4122) | Out-Null
$ws.Range("A215").PasteSpecial(-4122) | Out-Null

$ws.Range("E213").Copy() | Out-Null
$ws.Range("E214").PasteSpecial(-4122) | Out-Null
$ws.Range("E215").PasteSpecial(-4122) | Out-Null

$ws.Range("A214").Value = 212
$ws.Range("B214").Value = 7609150
$ws.Range("C214").Value = "Paraguay Division Profesional"
$ws.Range("D214").Value = "Paraguay Division Profesional"
$ws.Range("E214").Value = 45389.77083333334
$ws.Range("F214").Value = "Sol de America"
$ws.Range("G214").Value = "Olimpia Asuncion"
$ws.Range("K214").Value = 3.75
$ws.Range("L214").Value = 3.5
$ws.Range("M214").Value = 1.85
$ws.Range("N214").Value = 3.3
$ws.Range("O214").Value = 3.4
$ws.Range("P214").Value = 2
$ws.Range("Q214").Value = 0.25
$ws.Range("R214").Value = 2
$ws.Range("S214").Value = 1.8
$ws.Range("T214").Value = 2.25
$ws.Range("U214").Value = 1.825
$ws.Range("V214").Value = 1.975
$ws.Range("W214").Value = 0
$ws.Range("X214").Value = 0
$ws.Range("Y214").Value = 0
$ws.Range("Z214").Value = 0
$ws.Range("AA214").Value = 0

$ws.Range("A215").Value = 213
$ws.Range("B215").Value = 7609671
$ws.Range("C215").Value = "Paraguay Division Profesional"
$ws.Range("D215").Value = "Paraguay Division Profesional"
$ws.Range("E215").Value = 45389.875
$ws.Range("F215").Value = "2 de Mayo"
$ws.Range("G215").Value = "Sportivo Luqueno"
$ws.Range("K215").Value = 2.4
$ws.Range("L215").Value = 3.2
$ws.Range("M215").Value = 2.75
$ws.Range("N215").Value = 2.6
$ws.Range("O215").Value = 3.2
$ws.Range("P215").Value = 2.55
$ws.Range("Q215").Value = 0
$ws.Range("R215").Value = 1.925
$ws.Range("S215").Value = 1.875
$ws.Range("T215").Value = 2.25
$ws.Range("U215").Value = 1.95
$ws.Range("V215").Value = 1.85
$ws.Range("W215").Value = 0
$ws.Range("X215").Value = 0
$ws.Range("Y215").Value = 0
$ws.Range("Z215").Value = 0
$ws.Range("AA215").Value = 0
